# Commit: "Changed SummaryStats so that the last value has to be a
# detected result."
#
# On sheet "Alluvial for Mapping" (the SummaryStats table), columns X
# ("Last Cr") and Y ("Last Date") are updated for a number of wells so the
# reported "last" chromium result is a genuine detected value (not the
# 10.0 non-detect ceiling) -- with its matching sample date -- or, where no
# detected result exists at all, both cells become the literal text
# "No Detect Data".
#
# NOTE on technique: assigning a plain numeric- or date-looking string to
# Range.Value makes Excel auto-coerce it to a real number/date (and, for
# dates, silently mint a brand new date-formatted style for the cell).
# The source cells are plain General-formatted text cells, so to keep the
# values as literal text (and keep every cell on its original, shared
# style) we mark the cell as Text first, assign the value, and then
# restore the cell's original formatting by pasting the format (only)
# from an untouched cell that still carries the sheet's common style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alluvial for Mapping")

# Untouched cell that carries the same cell style ("s") as every cell we
# are about to edit -- used purely as a formatting donor so our edits
# don't leave the target cells on a newly minted (e.g. date-formatted)
# style.
$fmtSrc = $ws.Range("A1")

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $fmtSrc.Copy()
    $r.PasteSpecial(-4122)  # xlPasteFormats: restore the original style only
}

# Row 6 - CDBO-6
Set-TextValue "X6" "3.5"

# Row 7 - CDBO-7 (no detected result available)
Set-TextValue "X7" "No Detect Data"
Set-TextValue "Y7" "No Detect Data"

# Row 8 - MCA-1
Set-TextValue "X8" "2.8"
Set-TextValue "Y8" "2008-11-06"

# Row 15 - MCO-6B (no detected result available)
Set-TextValue "X15" "No Detect Data"
Set-TextValue "Y15" "No Detect Data"

# Row 21 - LAO-0.3 (no detected result available)
Set-TextValue "X21" "No Detect Data"
Set-TextValue "Y21" "No Detect Data"

# Row 22 - LAO-0.6
Set-TextValue "X22" "3.2"
Set-TextValue "Y22" "2008-08-29"

# Row 25 - LAO-1.6g
Set-TextValue "X25" "10.8"
Set-TextValue "Y25" "2008-08-27"

# Row 27 - LAO-4 (no detected result available)
Set-TextValue "X27" "No Detect Data"
Set-TextValue "Y27" "No Detect Data"

# Row 28 - LAO-4.5c
Set-TextValue "X28" "1.0"
Set-TextValue "Y28" "2007-07-19"

# Row 30 - LAO-B
Set-TextValue "X30" "1.3"
Set-TextValue "Y30" "2005-05-10"

# Row 31 - LAO-C
Set-TextValue "X31" "0.509"
Set-TextValue "Y31" "2004-06-03"

# Row 33 - PAO-4
Set-TextValue "X33" "2.6"
Set-TextValue "Y33" "2007-08-02"

# Row 36 - 18-MW-8
Set-TextValue "X36" "2.78"
Set-TextValue "Y36" "2009-09-08"

# Row 37 - 18-MW-9
Set-TextValue "X37" "1.5"
Set-TextValue "Y37" "2008-09-11"

# Row 38 - 3MAO-2
Set-TextValue "X38" "5.56"
Set-TextValue "Y38" "2010-06-07"

# Row 39 - PCAO-5 (no detected result available)
Set-TextValue "X39" "No Detect Data"
Set-TextValue "Y39" "No Detect Data"

# Row 41 - PCAO-7a (no detected result available)
Set-TextValue "X41" "No Detect Data"
Set-TextValue "Y41" "No Detect Data"

# Row 42 - PCAO-7b2
Set-TextValue "X42" "4.4"
Set-TextValue "Y42" "2009-03-06"

# Row 43 - PCAO-7c
Set-TextValue "X43" "2.74"
Set-TextValue "Y43" "2009-09-14"

# Row 44 - PCAO-9
Set-TextValue "X44" "9.1"
Set-TextValue "Y44" "2008-09-17"

# Row 46 - PCO-2
Set-TextValue "X46" "3.4"
Set-TextValue "Y46" "2007-12-06"

# Row 47 - PCO-3 (Y47 already holds the detected-result date; untouched)
Set-TextValue "X47" "2.6"

# Row 48 - TMO-1 (no detected result available)
Set-TextValue "X48" "No Detect Data"
Set-TextValue "Y48" "No Detect Data"

$excel.CutCopyMode = $false

Write-Host "SummaryStats: updated Last Cr / Last Date for 22 wells (44 cells)."
